$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.421.59'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.96%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.073.43'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.79%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '529.41'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +6.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.68'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +6.68%  '
$ws.Range("E8").Value = '  +5.42%  '
$ws.Range("E9").Value = '  +5.31%  '
$ws.Range("E10").Value = '  +7.80%  '
$ws.Range("E11").Value = '  +6.24%  '
$ws.Range("E12").Value = '  +2.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.595.17'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.49'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +9.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000175'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +17.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.281.09'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.79%  '
$ws.Range("E17").Value = '  +9.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.081.97'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.16%  '
$ws.Range("E19").Value = '  +6.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.22'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +5.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '343.41'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.09%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.75'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.508'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +8.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.60'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +5.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0₃0978'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +9.18%  '
$ws.Range("E27").Value = '  +4.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  +9.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.53'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +10.53%  '
$ws.Range("E31").Value = '  +7.20%  '
$ws.Range("E32").Value = '  +5.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.28'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +4.48%  '
$ws.Range("E34").Value = '  +9.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.43'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.01'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +7.57%  '
$ws.Range("E37").Value = '  +4.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.31'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +14.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0703'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.107.57'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.85%  '
$ws.Range("E41").Value = '  +4.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.98'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +12.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.671'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.00%  '
$ws.Range("E44").Value = '  +5.89%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.05'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +5.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.340.73'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +4.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.02'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.37%  '
$ws.Range("E49").Value = '  +5.84%  '
$ws.Range("E50").Value = '  +3.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.21'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +6.00%  '
